$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2605.5
$ws.Range("J17").Value = 2605.5
$ws.Range("L17").Value = 7816.5
$ws.Range("N17").Value = -8152.5
$ws.Range("H48").Value = 2943.6667
$ws.Range("I48").Value = 1817
$ws.Range("J48").Value = 3507
$ws.Range("K48").Value = 5451
$ws.Range("L48").Value = 10521
$ws.Range("M48").Value = -5159
$ws.Range("N48").Value = -11105
$ws.Range("H56").Value = 2943.6667
$ws.Range("I56").Value = 1817
$ws.Range("J56").Value = 3507
$ws.Range("K56").Value = 5451
$ws.Range("L56").Value = 10521
$ws.Range("M56").Value = -4917
$ws.Range("N56").Value = -11589
$ws.Range("H98").Value = 2686.7932
$ws.Range("I98").Value = 2823.8696
$ws.Range("J98").Value = 2161.3333
$ws.Range("K98").Value = 2823.8696
$ws.Range("L98").Value = 2161.3333
$ws.Range("M98").Value = -1325.8696
$ws.Range("N98").Value = -5157.3333
$ws.Range("H104").Value = 195.75
$ws.Range("I104").Value = 195.75
$ws.Range("K104").Value = 587.25
$ws.Range("M104").Value = 1159.75
$ws.Range("H112").Value = 30941.227
$ws.Range("J112").Value = 39618.332
$ws.Range("L112").Value = 118854.996
$ws.Range("N112").Value = -121070.996
$ws.Range("H122").Value = 2686.7932
$ws.Range("I122").Value = 2823.8696
$ws.Range("J122").Value = 2161.3333
$ws.Range("K122").Value = 8471.6088
$ws.Range("L122").Value = 6483.999899999999
$ws.Range("M122").Value = -6021.6088
$ws.Range("N122").Value = -11383.9999
$ws.Range("H129").Value = 48998.2
$ws.Range("I129").Value = 91479.2
$ws.Range("J129").Value = 6517.2
$ws.Range("K129").Value = 274437.6
$ws.Range("L129").Value = 19551.6
$ws.Range("M129").Value = -269437.6
$ws.Range("N129").Value = -29551.6
$ws.Range("H132").Value = 1917.9166
$ws.Range("I132").Value = 1626.3636
$ws.Range("K132").Value = 4879.0908
$ws.Range("M132").Value = -2349.0908
$ws.Range("H135").Value = 3935.4285
$ws.Range("J135").Value = 10666.667
$ws.Range("L135").Value = 96000.003
$ws.Range("N135").Value = -101070.003
$ws.Range("H137").Value = 1970.42
$ws.Range("I137").Value = 1282.9395
$ws.Range("J137").Value = 3304.9412
$ws.Range("K137").Value = 3848.8185
$ws.Range("L137").Value = 9914.8236
$ws.Range("M137").Value = -1298.8185
$ws.Range("N137").Value = -15014.8236
$ws.Range("H138").Value = 2576.3264
$ws.Range("I138").Value = 952.64
$ws.Range("J138").Value = 4267.6665
$ws.Range("K138").Value = 2857.92
$ws.Range("L138").Value = 12802.9995
$ws.Range("M138").Value = 2282.08
$ws.Range("N138").Value = -23082.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24933.516
$ws.Range("I32").Value = 22464.074
$ws.Range("K32").Value = 22464.074
$ws.Range("M32").Value = -22177.074
$ws.Range("H45").Value = 2182.8333
$ws.Range("I45").Value = 1212
$ws.Range("J45").Value = 4124.5
$ws.Range("K45").Value = 1212
$ws.Range("L45").Value = 4124.5
$ws.Range("M45").Value = -835
$ws.Range("N45").Value = -4878.5
$ws.Range("H74").Value = 102611
$ws.Range("I74").Value = 92303.5
$ws.Range("J74").Value = 115495.375
$ws.Range("K74").Value = 92303.5
$ws.Range("L74").Value = 115495.375
$ws.Range("M74").Value = -91429.5
$ws.Range("N74").Value = -117243.375
$ws.Range("H77").Value = 102611
$ws.Range("I77").Value = 92303.5
$ws.Range("J77").Value = 115495.375
$ws.Range("K77").Value = 461517.5
$ws.Range("L77").Value = 577476.875
$ws.Range("M77").Value = -457149.5
$ws.Range("N77").Value = -586212.875
$ws.Range("H132").Value = 23140.46
$ws.Range("I132").Value = 35251.375
$ws.Range("J132").Value = 3763
$ws.Range("K132").Value = 105754.125
$ws.Range("L132").Value = 11289
$ws.Range("M132").Value = -103224.125
$ws.Range("N132").Value = -16349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1655.9166
$ws.Range("I3").Value = 1928.2222
$ws.Range("J3").Value = 839
$ws.Range("K3").Value = 1928.2222
$ws.Range("L3").Value = 839
$ws.Range("M3").Value = -1814.2222
$ws.Range("N3").Value = -1067
$ws.Range("H20").Value = 2699.4211
$ws.Range("I20").Value = 2147.25
$ws.Range("J20").Value = 3646
$ws.Range("K20").Value = 2147.25
$ws.Range("L20").Value = 3646
$ws.Range("M20").Value = -1900.25
$ws.Range("N20").Value = -4140
$ws.Range("H99").Value = 4002
$ws.Range("J99").Value = 5499
$ws.Range("L99").Value = 5499
$ws.Range("N99").Value = -8495
$ws.Range("H134").Value = 25878.895
$ws.Range("I134").Value = 30959.322
$ws.Range("J134").Value = 3379.8572
$ws.Range("K134").Value = 92877.966
$ws.Range("L134").Value = 10139.5716
$ws.Range("M134").Value = -90342.966
$ws.Range("N134").Value = -15209.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2011.5209
$ws.Range("I31").Value = 1445.0278
$ws.Range("K31").Value = 1445.0278
$ws.Range("M31").Value = -1150.0278
$ws.Range("H34").Value = 2011.5209
$ws.Range("I34").Value = 1445.0278
$ws.Range("K34").Value = 1445.0278
$ws.Range("M34").Value = -1243.0278
$ws.Range("H58").Value = 7855.5386
$ws.Range("I58").Value = 7414.1665
$ws.Range("J58").Value = 8848.625
$ws.Range("K58").Value = 7414.1665
$ws.Range("L58").Value = 8848.625
$ws.Range("M58").Value = -7211.1665
$ws.Range("N58").Value = -9254.625
$ws.Range("H132").Value = 2518.4
$ws.Range("I132").Value = 2666.6843
$ws.Range("J132").Value = 2048.8333
$ws.Range("K132").Value = 8000.0529
$ws.Range("L132").Value = 6146.499899999999
$ws.Range("M132").Value = -5470.0529
$ws.Range("N132").Value = -11206.4999
$ws.Range("H134").Value = 1894.0869
$ws.Range("I134").Value = 1133.9524
$ws.Range("K134").Value = 3401.857199999999
$ws.Range("M134").Value = -866.8571999999995
$ws.Range("H136").Value = 7855.5386
$ws.Range("I136").Value = 7414.1665
$ws.Range("J136").Value = 8848.625
$ws.Range("K136").Value = 22242.4995
$ws.Range("L136").Value = 26545.875
$ws.Range("M136").Value = -19692.4995
$ws.Range("N136").Value = -31645.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3651
$ws.Range("I134").Value = 3501.111
$ws.Range("K134").Value = 10503.333
$ws.Range("M134").Value = -5433.332999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20757.223
$ws.Range("I43").Value = 3908
$ws.Range("J43").Value = 25571.285
$ws.Range("K43").Value = 3908
$ws.Range("L43").Value = 25571.285
$ws.Range("M43").Value = -3757
$ws.Range("N43").Value = -25873.285
$ws.Range("H63").Value = 16862.25
$ws.Range("I63").Value = 15898
$ws.Range("K63").Value = 15898
$ws.Range("M63").Value = -15212
$ws.Range("H66").Value = 16862.25
$ws.Range("I66").Value = 15898
$ws.Range("K66").Value = 47694
$ws.Range("M66").Value = -44262

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49952.844
$ws.Range("I7").Value = 55358.824
$ws.Range("K7").Value = 55358.824
$ws.Range("M7").Value = -55246.824
$ws.Range("H46").Value = 2643.959
$ws.Range("J46").Value = 2807
$ws.Range("L46").Value = 2807
$ws.Range("N46").Value = -3183
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -20956
$ws.Range("H61").Value = 4895.25
$ws.Range("I61").Value = 4271.375
$ws.Range("J61").Value = 6143
$ws.Range("K61").Value = 4271.375
$ws.Range("L61").Value = 6143
$ws.Range("M61").Value = -4069.375
$ws.Range("N61").Value = -6547
$ws.Range("H100").Value = 5626
$ws.Range("I100").Value = 3501
$ws.Range("K100").Value = 3501
$ws.Range("M100").Value = -2960
$ws.Range("H113").Value = 4895.25
$ws.Range("I113").Value = 4271.375
$ws.Range("J113").Value = 6143
$ws.Range("K113").Value = 4271.375
$ws.Range("L113").Value = 6143
$ws.Range("M113").Value = -2101.375
$ws.Range("N113").Value = -10483
$ws.Range("H122").Value = 3793.8845
$ws.Range("I122").Value = 3550.6
$ws.Range("J122").Value = 9876
$ws.Range("K122").Value = 10651.8
$ws.Range("L122").Value = 29628
$ws.Range("M122").Value = -8201.799999999999
$ws.Range("N122").Value = -34528
$ws.Range("H126").Value = 49952.844
$ws.Range("I126").Value = 55358.824
$ws.Range("K126").Value = 166076.472
$ws.Range("M126").Value = -163606.472
$ws.Range("H132").Value = 3233.2424
$ws.Range("I132").Value = 2444.625
$ws.Range("K132").Value = 7333.875
$ws.Range("M132").Value = -4803.875
$ws.Range("H136").Value = 3848.5417
$ws.Range("I136").Value = 2302.3157
$ws.Range("K136").Value = 6906.9471
$ws.Range("M136").Value = -4356.9471

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 405517.8
$ws.Range("I62").Value = 669196.7
$ws.Range("J62").Value = 9999.5
$ws.Range("K62").Value = 669196.7
$ws.Range("L62").Value = 9999.5
$ws.Range("M62").Value = -668572.7
$ws.Range("N62").Value = -11247.5
$ws.Range("H65").Value = 405517.8
$ws.Range("I65").Value = 669196.7
$ws.Range("J65").Value = 9999.5
$ws.Range("K65").Value = 3345983.5
$ws.Range("L65").Value = 49997.5
$ws.Range("M65").Value = -3342863.5
$ws.Range("N65").Value = -56237.5
$ws.Range("H96").Value = 8850
$ws.Range("J96").Value = 1500
$ws.Range("L96").Value = 1500
$ws.Range("N96").Value = -4246
$ws.Range("H107").Value = 13122
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H132").Value = 5901.8184
$ws.Range("I132").Value = 4726.2144
$ws.Range("K132").Value = 14178.6432
$ws.Range("M132").Value = -11648.6432
